# Budget.xlsx edits:
# - H13: price 13.26 -> 13.28
# - H22: price 18.81 -> 12.55
# - H25/I25: clear the price values (no longer priced individually)
# - H29: price 8.46 -> 12.22
# - leave the selection on H24, matching the last cell touched

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H13").Value = 13.28

$ws.Range("H22").Value = 12.55

$ws.Range("H25:I25").ClearContents()

$ws.Range("H29").Value = 12.22

$ws.Range("H24").Select() | Out-Null
